$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.194.61"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "2.056.37"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'248.46"
$ws.Range("E5").Value = "  -1.81%  "

$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'57.09"
$ws.Range("E8").Value = "  -2.38%  "

$ws.Range("D9").Value = "'0.386"
$ws.Range("E9").Value = "  -0.14%  "

$ws.Range("D10").Value = "'0.0786"
$ws.Range("E10").Value = "  -1.52%  "

$ws.Range("D12").Value = "'16.32"
$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("D14").Value = "2.353.62"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").Value = "'5.80"
$ws.Range("E15").Value = "  +3.11%  "

$ws.Range("D16").Value = "2.055.94"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").Value = "'18.78"
$ws.Range("E17").Value = "  +13.22%  "

$ws.Range("D18").Value = "37.239.06"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("D19").Value = "'74.88"
$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  -1.79%  "

$ws.Range("D21").Value = "'5.50"
$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("D22").Value = "'237.80"
$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'2.49"
$ws.Range("E24").Value = "  +4.08%  "

$ws.Range("E25").Value = "  +4.07%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'170.56"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.18"
$ws.Range("E27").Value = "  -5.54%  "

$ws.Range("D28").Value = "'20.26"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("E29").Value = "  -0.94%  "

$ws.Range("D30").Value = "'5.23"
$ws.Range("E30").Value = "  +9.45%  "

$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +2.39%  "

$ws.Range("D32").Value = "'0.0626"
$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("D33").Value = "'4.60"
$ws.Range("E33").Value = "  +2.07%  "

$ws.Range("D34").Value = "'0.0887"
$ws.Range("E34").Value = "  +0.27%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").Value = "'2.28"

$ws.Range("E37").Value = "  +0.99%  "

$ws.Range("E38").Value = "  -2.00%  "

$ws.Range("D39").Value = "'5.30"
$ws.Range("E39").Value = "  +15.85%  "

$ws.Range("E40").Value = "  +7.68%  "

$ws.Range("D41").Value = "'0.0999"
$ws.Range("E41").Value = "  -12.78%  "

$ws.Range("E42").Value = "  -0.42%  "

$ws.Range("E43").Value = "  -0.68%  "

$ws.Range("D44").Value = "'1.16"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").Value = "'96.80"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("E46").Value = "  -1.33%  "

$ws.Range("D47").Value = "1.277.05"
$ws.Range("E47").Value = "  -1.28%  "

$ws.Range("E48").Value = "  -2.12%  "

$ws.Range("D49").Value = "'6.87"
$ws.Range("E49").Value = "  -0.75%  "

$ws.Range("D50").Value = "2.239.99"
$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("D51").Value = "'44.42"
$ws.Range("E51").Value = "  +0.70%  "
